{"js": "// Add a new closing title paragraph \"Cocumento Game Design\" at the very\n// end of the document body, right after the last (empty) paragraph and\n// before the section break \u2014 matching the author's commit\n// \"T\u00edtulo Documento Game Design\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"Cocumento Game Design\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Add a new closing title paragraph \"Cocumento Game Design\" at the very\n# end of the document, right after the last (empty) paragraph and before\n# the section break -- matching the author's commit\n# \"T\u00edtulo Documento Game Design\".\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Collapse(0)  # wdCollapseEnd\n$rng.InsertParagraphAfter()\n$rng.Collapse(0)\n$rng.InsertAfter(\"Cocumento Game Design\")\n"}
